$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 83.33333333333334
$ws.Range("C2").Value = 83.33333333333334
$ws.Range("D2").Value = 83.33333333333334
$ws.Range("E2").Value = 73.69791666666667
$ws.Range("F2").Value = 59.89583333333333
$ws.Range("G2").Value = 59.63541666666667
$ws.Range("H2").Value = 56.25
$ws.Range("I2").Value = 54.16666666666667
$ws.Range("J2").Value = 54.16666666666667
$ws.Range("K2").Value = 54.16666666666667
$ws.Range("L2").Value = 53.38541666666667
$ws.Range("M2").Value = 53.38541666666667
$ws.Range("N2").Value = 53.38541666666667
$ws.Range("O2").Value = 53.38541666666667
$ws.Range("P2").Value = 53.38541666666667
$ws.Range("Q2").Value = 53.38541666666667
$ws.Range("R2").Value = 53.38541666666667
$ws.Range("S2").Value = 53.38541666666667
$ws.Range("T2").Value = 53.38541666666667
$ws.Range("U2").Value = 53.38541666666667
